# Update the "master" output of the two-digit x two-digit multiplication
# worksheet: bump the dated heading by one day and regenerate all 25
# practice problems/answers in the 5x5 table.
$d = $word.ActiveDocument

# Heading date
$d.Content.Find.Execute("2024-03-16 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-17 Sunday", 2) | Out-Null

# Table cells (each "a×b=c" string is unique in the document, so a plain
# whole-string Find/Replace safely targets the correct run each time).
$d.Content.Find.Execute("39×93=3627", $true, $false, $false, $false, $false, $true, 1, $false, "28×83=2324", 2) | Out-Null
$d.Content.Find.Execute("83×71=5893", $true, $false, $false, $false, $false, $true, 1, $false, "96×39=3744", 2) | Out-Null
$d.Content.Find.Execute("77×84=6468", $true, $false, $false, $false, $false, $true, 1, $false, "62×36=2232", 2) | Out-Null
$d.Content.Find.Execute("38×37=1406", $true, $false, $false, $false, $false, $true, 1, $false, "22×32=704", 2) | Out-Null
$d.Content.Find.Execute("34×16=544", $true, $false, $false, $false, $false, $true, 1, $false, "94×41=3854", 2) | Out-Null
$d.Content.Find.Execute("98×64=6272", $true, $false, $false, $false, $false, $true, 1, $false, "52×97=5044", 2) | Out-Null
$d.Content.Find.Execute("38×52=1976", $true, $false, $false, $false, $false, $true, 1, $false, "39×84=3276", 2) | Out-Null
$d.Content.Find.Execute("51×70=3570", $true, $false, $false, $false, $false, $true, 1, $false, "68×21=1428", 2) | Out-Null
$d.Content.Find.Execute("60×54=3240", $true, $false, $false, $false, $false, $true, 1, $false, "96×14=1344", 2) | Out-Null
$d.Content.Find.Execute("46×75=3450", $true, $false, $false, $false, $false, $true, 1, $false, "53×45=2385", 2) | Out-Null
$d.Content.Find.Execute("88×81=7128", $true, $false, $false, $false, $false, $true, 1, $false, "26×12=312", 2) | Out-Null
$d.Content.Find.Execute("72×79=5688", $true, $false, $false, $false, $false, $true, 1, $false, "49×84=4116", 2) | Out-Null
$d.Content.Find.Execute("85×64=5440", $true, $false, $false, $false, $false, $true, 1, $false, "71×96=6816", 2) | Out-Null
$d.Content.Find.Execute("57×46=2622", $true, $false, $false, $false, $false, $true, 1, $false, "18×40=720", 2) | Out-Null
$d.Content.Find.Execute("55×77=4235", $true, $false, $false, $false, $false, $true, 1, $false, "24×65=1560", 2) | Out-Null
$d.Content.Find.Execute("74×54=3996", $true, $false, $false, $false, $false, $true, 1, $false, "17×37=629", 2) | Out-Null
$d.Content.Find.Execute("49×90=4410", $true, $false, $false, $false, $false, $true, 1, $false, "12×22=264", 2) | Out-Null
$d.Content.Find.Execute("60×69=4140", $true, $false, $false, $false, $false, $true, 1, $false, "69×84=5796", 2) | Out-Null
$d.Content.Find.Execute("85×11=935", $true, $false, $false, $false, $false, $true, 1, $false, "58×60=3480", 2) | Out-Null
$d.Content.Find.Execute("59×83=4897", $true, $false, $false, $false, $false, $true, 1, $false, "38×99=3762", 2) | Out-Null
$d.Content.Find.Execute("43×71=3053", $true, $false, $false, $false, $false, $true, 1, $false, "39×34=1326", 2) | Out-Null
$d.Content.Find.Execute("97×95=9215", $true, $false, $false, $false, $false, $true, 1, $false, "29×71=2059", 2) | Out-Null
$d.Content.Find.Execute("75×42=3150", $true, $false, $false, $false, $false, $true, 1, $false, "12×54=648", 2) | Out-Null
$d.Content.Find.Execute("60×46=2760", $true, $false, $false, $false, $false, $true, 1, $false, "47×87=4089", 2) | Out-Null
$d.Content.Find.Execute("60×12=720", $true, $false, $false, $false, $false, $true, 1, $false, "27×57=1539", 2) | Out-Null
